# Auto-generated edit script applying numeric corrections to the
# Anima_Profits leve-profit tables (ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets).
# Each block updates currentAveragePrice* / Leve* columns (H:N) for one row,
# matching the scheduled-runner price refresh described in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 55
$ws.Cells.Item(55, 8).Value = 45455400
$ws.Cells.Item(55, 9).Value = 71429704
$ws.Cells.Item(55, 11).Value = 71429704
$ws.Cells.Item(55, 13).Value = -71429490

# ALC row 62
$ws.Cells.Item(62, 8).Value = 7604.923
$ws.Cells.Item(62, 9).Value = 2923.5715
$ws.Cells.Item(62, 11).Value = 2923.5715
$ws.Cells.Item(62, 13).Value = -2299.5715

# ALC row 65
$ws.Cells.Item(65, 8).Value = 7604.923
$ws.Cells.Item(65, 9).Value = 2923.5715
$ws.Cells.Item(65, 11).Value = 14617.8575
$ws.Cells.Item(65, 13).Value = -11497.8575

# ALC row 132
$ws.Cells.Item(132, 8).Value = 7307.3335
$ws.Cells.Item(132, 9).Value = 7638.9
$ws.Cells.Item(132, 11).Value = 22916.7
$ws.Cells.Item(132, 13).Value = -20386.7

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61
$ws.Cells.Item(61, 8).Value = 8775967
$ws.Cells.Item(61, 9).Value = 18520952
$ws.Cells.Item(61, 10).Value = 5479.9
$ws.Cells.Item(61, 11).Value = 18520952
$ws.Cells.Item(61, 12).Value = 5479.9
$ws.Cells.Item(61, 13).Value = -18520740
$ws.Cells.Item(61, 14).Value = -5903.9

# ARM row 132
$ws.Cells.Item(132, 8).Value = 1835110
$ws.Cells.Item(132, 9).Value = 2676.5312
$ws.Cells.Item(132, 11).Value = 8029.5936
$ws.Cells.Item(132, 13).Value = -5499.5936

# ARM row 136
$ws.Cells.Item(136, 8).Value = 8775967
$ws.Cells.Item(136, 9).Value = 18520952
$ws.Cells.Item(136, 10).Value = 5479.9
$ws.Cells.Item(136, 11).Value = 55562856
$ws.Cells.Item(136, 12).Value = 16439.7
$ws.Cells.Item(136, 13).Value = -55560306
$ws.Cells.Item(136, 14).Value = -21539.7

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Cells.Item(31, 8).Value = 4701.7866
$ws.Cells.Item(31, 9).Value = 1123.4286
$ws.Cells.Item(31, 10).Value = 7832.85
$ws.Cells.Item(31, 11).Value = 1123.4286
$ws.Cells.Item(31, 12).Value = 7832.85
$ws.Cells.Item(31, 13).Value = -828.4286
$ws.Cells.Item(31, 14).Value = -8422.85

# CRP row 34
$ws.Cells.Item(34, 8).Value = 4701.7866
$ws.Cells.Item(34, 9).Value = 1123.4286
$ws.Cells.Item(34, 10).Value = 7832.85
$ws.Cells.Item(34, 11).Value = 1123.4286
$ws.Cells.Item(34, 12).Value = 7832.85
$ws.Cells.Item(34, 13).Value = -921.4286
$ws.Cells.Item(34, 14).Value = -8236.85

# CRP row 36
$ws.Cells.Item(36, 8).Value = 2774
$ws.Cells.Item(36, 9).Value = 2774
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 2774
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -2386
$ws.Cells.Item(36, 14).Value = ""

# CRP row 40
$ws.Cells.Item(40, 8).Value = 2774
$ws.Cells.Item(40, 9).Value = 2774
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 2774
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -2614
$ws.Cells.Item(40, 14).Value = ""

# CRP row 54
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 13).Value = ""
$ws.Cells.Item(54, 14).Value = ""

# CRP row 134
$ws.Cells.Item(134, 8).Value = 8934447
$ws.Cells.Item(134, 9).Value = 13164554
$ws.Cells.Item(134, 11).Value = 39493662
$ws.Cells.Item(134, 13).Value = -39491127

# CRP row 141
$ws.Cells.Item(141, 8).Value = 241664.44
$ws.Cells.Item(141, 9).Value = 200000
$ws.Cells.Item(141, 10).Value = 246872.5
$ws.Cells.Item(141, 11).Value = 200000
$ws.Cells.Item(141, 12).Value = 246872.5
$ws.Cells.Item(141, 13).Value = -194820
$ws.Cells.Item(141, 14).Value = -257232.5

$ws = $wb.Worksheets.Item("CUL")
# CUL row 81
$ws.Cells.Item(81, 8).Value = 7519.154
$ws.Cells.Item(81, 9).Value = 1612.5
$ws.Cells.Item(81, 11).Value = 4837.5
$ws.Cells.Item(81, 13).Value = -3714.5

# CUL row 84
$ws.Cells.Item(84, 8).Value = 7519.154
$ws.Cells.Item(84, 9).Value = 1612.5
$ws.Cells.Item(84, 11).Value = 14512.5
$ws.Cells.Item(84, 13).Value = -8896.5

# CUL row 92
$ws.Cells.Item(92, 8).Value = 1000
$ws.Cells.Item(92, 10).Value = 1000
$ws.Cells.Item(92, 12).Value = 3000
$ws.Cells.Item(92, 14).Value = -5496

# CUL row 131
$ws.Cells.Item(131, 8).Value = 4252.0293
$ws.Cells.Item(131, 10).Value = 6322.727
$ws.Cells.Item(131, 12).Value = 18968.181
$ws.Cells.Item(131, 14).Value = -29048.181

$ws = $wb.Worksheets.Item("GSM")
# GSM row 2
$ws.Cells.Item(2, 8).Value = 262.125
$ws.Cells.Item(2, 9).Value = 271
$ws.Cells.Item(2, 10).Value = 200
$ws.Cells.Item(2, 11).Value = 271
$ws.Cells.Item(2, 12).Value = 200
$ws.Cells.Item(2, 13).Value = -158
$ws.Cells.Item(2, 14).Value = -426

# GSM row 18
$ws.Cells.Item(18, 8).Value = 8300
$ws.Cells.Item(18, 10).Value = 8300
$ws.Cells.Item(18, 12).Value = 8300
$ws.Cells.Item(18, 14).Value = -8886

# GSM row 26
$ws.Cells.Item(26, 8).Value = 31410.166
$ws.Cells.Item(26, 9).Value = 7999.75
$ws.Cells.Item(26, 10).Value = 78231
$ws.Cells.Item(26, 11).Value = 7999.75
$ws.Cells.Item(26, 12).Value = 78231
$ws.Cells.Item(26, 13).Value = -7719.75
$ws.Cells.Item(26, 14).Value = -78791

# GSM row 43
$ws.Cells.Item(43, 8).Value = 5729.4443
$ws.Cells.Item(43, 9).Value = 1094.6666
$ws.Cells.Item(43, 10).Value = 14999
$ws.Cells.Item(43, 11).Value = 1094.6666
$ws.Cells.Item(43, 12).Value = 14999
$ws.Cells.Item(43, 13).Value = -943.6666
$ws.Cells.Item(43, 14).Value = -15301

# GSM row 46
$ws.Cells.Item(46, 8).Value = 4499.9443
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 4499.9443
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 4499.9443
$ws.Cells.Item(46, 13).Value = ""
$ws.Cells.Item(46, 14).Value = -4811.9443

# GSM row 50
$ws.Cells.Item(50, 8).Value = 31410.166
$ws.Cells.Item(50, 9).Value = 7999.75
$ws.Cells.Item(50, 10).Value = 78231
$ws.Cells.Item(50, 11).Value = 7999.75
$ws.Cells.Item(50, 12).Value = 78231
$ws.Cells.Item(50, 13).Value = -7501.75
$ws.Cells.Item(50, 14).Value = -79227

# GSM row 102
$ws.Cells.Item(102, 8).Value = 1196.1875
$ws.Cells.Item(102, 9).Value = 1173.2142
$ws.Cells.Item(102, 10).Value = 1357
$ws.Cells.Item(102, 11).Value = 1173.2142
$ws.Cells.Item(102, 12).Value = 1357
$ws.Cells.Item(102, 13).Value = 448.7858000000001
$ws.Cells.Item(102, 14).Value = -4601

# GSM row 132
$ws.Cells.Item(132, 8).Value = 2751.5518
$ws.Cells.Item(132, 9).Value = 1963.8572
$ws.Cells.Item(132, 10).Value = 4819.25
$ws.Cells.Item(132, 11).Value = 5891.571599999999
$ws.Cells.Item(132, 12).Value = 14457.75
$ws.Cells.Item(132, 13).Value = -3361.571599999999
$ws.Cells.Item(132, 14).Value = -19517.75

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Cells.Item(22, 8).Value = 12530.134
$ws.Cells.Item(22, 9).Value = 650
$ws.Cells.Item(22, 10).Value = 18470.2
$ws.Cells.Item(22, 11).Value = 650
$ws.Cells.Item(22, 12).Value = 18470.2
$ws.Cells.Item(22, 13).Value = -355
$ws.Cells.Item(22, 14).Value = -19060.2

# LTW row 27
$ws.Cells.Item(27, 8).Value = 12530.134
$ws.Cells.Item(27, 9).Value = 650
$ws.Cells.Item(27, 10).Value = 18470.2
$ws.Cells.Item(27, 11).Value = 650
$ws.Cells.Item(27, 12).Value = 18470.2
$ws.Cells.Item(27, 13).Value = -543
$ws.Cells.Item(27, 14).Value = -18684.2

$ws = $wb.Worksheets.Item("WVR")
# WVR row 139
$ws.Cells.Item(139, 8).Value = 62265.176
$ws.Cells.Item(139, 9).Value = 61181
$ws.Cells.Item(139, 10).Value = 62332.938
$ws.Cells.Item(139, 11).Value = 61181
$ws.Cells.Item(139, 12).Value = 62332.938
$ws.Cells.Item(139, 13).Value = -56041
$ws.Cells.Item(139, 14).Value = -72612.93799999999

